# Include the image title (from the markdown `![alt](link "title")`) in
# PowerPoint's description of the image, alongside the existing alt text.
#
# The picture "Picture 1" on slide 1 currently has AlternativeText
# "lalune.jpg" (the link); we append the title ("fig:") so the
# description becomes "fig:  lalune.jpg".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Picture 1")
$shape.AlternativeText = "fig:  lalune.jpg"
